$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = '2026-02-17 23:48:39'
$ws.Cells.Item(2, 14).Value = '-0.3 °C 23:19 TU'
$ws.Cells.Item(2, 15).Value = '2.4 °C'
$ws.Cells.Item(3, 5).Value = '2026-02-17 23:48:41'
$ws.Cells.Item(4, 5).Value = '2026-02-17 23:48:44'
$ws.Cells.Item(4, 10).Value = '1018.8 hPa'
$ws.Cells.Item(5, 5).Value = '2026-02-17 23:48:47'
$ws.Cells.Item(6, 5).Value = '2026-02-17 23:48:50'
$ws.Cells.Item(6, 10).Value = '1018.7 hPa'
$ws.Cells.Item(7, 5).Value = '2026-02-17 23:48:52'
$ws.Cells.Item(7, 10).Value = '1018.6 hPa'
$ws.Cells.Item(8, 5).Value = '2026-02-17 23:48:55'
$ws.Cells.Item(9, 5).Value = '2026-02-17 23:48:58'
$ws.Cells.Item(9, 8).Value = '63%'
$ws.Cells.Item(10, 5).Value = '2026-02-17 23:49:00'
$ws.Cells.Item(11, 5).Value = '2026-02-17 23:49:03'
$ws.Cells.Item(11, 8).Value = '56%'
$ws.Cells.Item(11, 15).Value = '6.9 °C'
$ws.Cells.Item(12, 5).Value = '2026-02-17 23:49:06'
$ws.Cells.Item(12, 8).Value = '65%'
$ws.Cells.Item(12, 15).Value = '12.2 °C'
$ws.Cells.Item(13, 5).Value = '2026-02-17 23:49:08'
$ws.Cells.Item(13, 8).Value = '49%'
$ws.Cells.Item(13, 10).Value = '1018.4 hPa'
$ws.Cells.Item(13, 14).Value = '0.7 °C 23:29 TU'
$ws.Cells.Item(13, 15).Value = '6.4 °C'
$ws.Cells.Item(14, 5).Value = '2026-02-17 23:49:10'
$ws.Cells.Item(14, 8).Value = '72%'
$ws.Cells.Item(15, 5).Value = '2026-02-17 23:49:13'
$ws.Cells.Item(15, 8).Value = '62%'
$ws.Cells.Item(15, 15).Value = '11.6 °C'
$ws.Cells.Item(16, 5).Value = '2026-02-17 23:49:16'
$ws.Cells.Item(16, 15).Value = '-3.0 °C'
$ws.Cells.Item(17, 5).Value = '2026-02-17 23:49:19'
$ws.Cells.Item(17, 8).Value = '84%'
$ws.Cells.Item(18, 5).Value = '2026-02-17 23:49:21'
$ws.Cells.Item(19, 5).Value = '2026-02-17 23:49:24'
$ws.Cells.Item(19, 8).Value = '77%'
$ws.Cells.Item(19, 15).Value = '7.1 °C'
$ws.Cells.Item(20, 5).Value = '2026-02-17 23:49:27'
$ws.Cells.Item(20, 8).Value = '69%'
$ws.Cells.Item(21, 5).Value = '2026-02-17 23:49:29'
$ws.Cells.Item(21, 10).Value = '1017.4 hPa'
$ws.Cells.Item(21, 14).Value = '3.8 °C 23:28 TU'
$ws.Cells.Item(21, 15).Value = '9.2 °C'
$ws.Cells.Item(22, 5).Value = '2026-02-17 23:49:32'
$ws.Cells.Item(23, 5).Value = '2026-02-17 23:49:35'
$ws.Cells.Item(23, 15).Value = '-3.4 °C'
$ws.Cells.Item(24, 5).Value = '2026-02-17 23:49:37'
$ws.Cells.Item(24, 10).Value = '1018.8 hPa'
$ws.Cells.Item(24, 15).Value = '12.5 °C'
$ws.Cells.Item(25, 5).Value = '2026-02-17 23:49:40'
$ws.Cells.Item(26, 5).Value = '2026-02-17 23:49:43'
$ws.Cells.Item(27, 5).Value = '2026-02-17 23:49:45'
$ws.Cells.Item(27, 8).Value = '57%'
$ws.Cells.Item(27, 15).Value = '-0.3 °C'
$ws.Cells.Item(28, 5).Value = '2026-02-17 23:49:48'
$ws.Cells.Item(28, 10).Value = '1018.6 hPa'
$ws.Cells.Item(29, 5).Value = '2026-02-17 23:49:50'
$ws.Cells.Item(29, 15).Value = '11.7 °C'
$ws.Cells.Item(30, 5).Value = '2026-02-17 23:49:53'
$ws.Cells.Item(30, 8).Value = '68%'
$ws.Cells.Item(30, 14).Value = '7.0 °C 23:22 TU'
$ws.Cells.Item(30, 15).Value = '10.9 °C'
$ws.Cells.Item(31, 5).Value = '2026-02-17 23:49:56'
$ws.Cells.Item(31, 8).Value = '69%'
$ws.Cells.Item(32, 5).Value = '2026-02-17 23:49:58'
$ws.Cells.Item(33, 5).Value = '2026-02-17 23:50:01'
$ws.Cells.Item(33, 8).Value = '47%'
$ws.Cells.Item(33, 10).Value = '1017.8 hPa'
$ws.Cells.Item(33, 14).Value = '1.9 °C 23:12 TU'
$ws.Cells.Item(33, 15).Value = '6.1 °C'
$ws.Cells.Item(34, 5).Value = '2026-02-17 23:50:03'
$ws.Cells.Item(35, 5).Value = '2026-02-17 23:50:06'
$ws.Cells.Item(35, 8).Value = '92%'
$ws.Cells.Item(36, 5).Value = '2026-02-17 23:50:08'
$ws.Cells.Item(36, 8).Value = '63%'
$ws.Cells.Item(37, 5).Value = '2026-02-17 23:50:11'
$ws.Cells.Item(37, 8).Value = '75%'
$ws.Cells.Item(37, 15).Value = '7.0 °C'
$ws.Cells.Item(38, 5).Value = '2026-02-17 23:50:14'
$ws.Cells.Item(38, 8).Value = '79%'
$ws.Cells.Item(39, 5).Value = '2026-02-17 23:50:16'
$ws.Cells.Item(39, 13).Value = '0.4 °C 23:29 TU'
$ws.Cells.Item(39, 15).Value = '-2.4 °C'
$ws.Cells.Item(40, 5).Value = '2026-02-17 23:50:19'
$ws.Cells.Item(40, 8).Value = '55%'
$ws.Cells.Item(40, 10).Value = '1018.3 hPa'
$ws.Cells.Item(40, 15).Value = '9.1 °C'
$ws.Cells.Item(41, 5).Value = '2026-02-17 23:50:21'
$ws.Cells.Item(41, 8).Value = '54%'
$ws.Cells.Item(41, 10).Value = '1018.3 hPa'
$ws.Cells.Item(41, 14).Value = '11.2 °C 23:29 TU'
$ws.Cells.Item(41, 15).Value = '16.0 °C'
$ws.Cells.Item(42, 5).Value = '2026-02-17 23:50:23'
$ws.Cells.Item(42, 15).Value = '12.5 °C'
$ws.Cells.Item(43, 5).Value = '2026-02-17 23:50:26'
$ws.Cells.Item(44, 5).Value = '2026-02-17 23:50:29'
$ws.Cells.Item(44, 8).Value = '78%'
$ws.Cells.Item(45, 5).Value = '2026-02-17 23:50:31'
$ws.Cells.Item(45, 14).Value = '1.4 °C 23:29 TU'
$ws.Cells.Item(45, 15).Value = '5.0 °C'
$ws.Cells.Item(46, 5).Value = '2026-02-17 23:50:34'
$ws.Cells.Item(46, 8).Value = '61%'
$ws.Cells.Item(46, 14).Value = '8.7 °C 23:28 TU'
$ws.Cells.Item(46, 15).Value = '14.8 °C'
